$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22 (22-Aug-23): everyone PRESENT ---
$ws.Range("A22").Value = 45160
$ws.Range("A22").NumberFormat = $ws.Range("A21").NumberFormat

$row22 = @("PRESENT","PRESENT","PRESENT","PRESENT","PRESENT","PRESENT","PRESENT","PRESENT","PRESENT","PRESENT")
$cols = @("B","C","D","E","F","G","H","I","J","K")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "22").Value = $row22[$i]
}

# --- Row 23 (23-Aug-23): mixed attendance ---
$ws.Range("A23").Value = 45161
$ws.Range("A23").NumberFormat = $ws.Range("A21").NumberFormat

$row23 = @{
    "B" = "PRESENT"
    "C" = "PRESENT"
    "D" = "ABSENT"
    "E" = "PRESENT"
    "F" = "PRESENT"
    "G" = "ABSENT"
    "H" = "ABSENT"
    "I" = "ABSENT"
    "J" = "ABSENT"
    "K" = "ABSENT"
}
foreach ($col in $cols) {
    $ws.Range($col + "23").Value = $row23[$col]
}

# --- Comments explaining the ABSENT marks on 23-Aug-23 ---
$ws.Range("D23").AddComment("LENOVO:" + [char]10 + "outside" + [char]10)
$ws.Range("G23").AddComment("LENOVO:" + [char]10 + "Headeche")
$ws.Range("H23").AddComment("LENOVO:" + [char]10 + "No response")
$ws.Range("I23").AddComment("LENOVO:" + [char]10 + "No response")
$ws.Range("J23").AddComment("LENOVO:" + [char]10 + "No response")
$ws.Range("K23").AddComment("LENOVO:" + [char]10 + "having some work")

# --- Update the view: scroll the frozen pane down and move the selection ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $true
$ws.Range("G28").Select()
